# Sync attendance_reports: reorder the "Recorded By" (column G) names so
# that System-type entries are listed before other (email) entries.
#
# Sort rule observed in the target diff:
#   - Split the comma-separated list of names in each G cell.
#   - "system"/"System" entries sort first; among themselves they are
#     ordered case-sensitively (so "System" precedes "system").
#   - All remaining (email) entries follow, ordered case-insensitively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function IsSystemName($p) {
    return ($p.ToLower() -eq "system")
}

function CompareNames($a, $b) {
    $aSys = IsSystemName $a
    $bSys = IsSystemName $b
    if ($aSys -and -not $bSys) { return -1 }
    if ((-not $aSys) -and $bSys) { return 1 }
    if ($aSys -and $bSys) {
        return $a.CompareTo($b)
    }
    return $a.ToLower().CompareTo($b.ToLower())
}

function SortNames($s) {
    $parts = $s -split ","
    $list = New-Object System.Collections.ArrayList
    foreach ($p in $parts) { [void]$list.Add($p.Trim()) }

    # insertion sort, stable, using CompareNames
    for ($i = 1; $i -lt $list.Count; $i++) {
        $key = $list[$i]
        $j = $i - 1
        while ($j -ge 0 -and (CompareNames $list[$j] $key) -gt 0) {
            $list[$j + 1] = $list[$j]
            $j = $j - 1
        }
        $list[$j + 1] = $key
    }
    return ($list -join ", ")
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Value2
    # NOTE: the -eq/-ne comparison operators in this host are
    # case-insensitive, so use the case-sensitive .Equals() method to
    # decide whether a real (case-sensitive) change is needed.
    if ($value -ne $null -and $value -ne "" -and $value.Contains(",")) {
        $newValue = SortNames $value
        if (-not $newValue.Equals($value)) {
            $cell.Value = $newValue
        }
    }
}
